# Fix cost bug: update constant-value rows (columns J:AS) on the single
# worksheet "strategy_id-0" for several prodinit_ippu_* variable rows.
#
# Each affected row currently holds the same old constant value repeated
# across columns J through AS; it must be replaced by a new constant
# value repeated across the same columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowUpdates = @{
    96  = 33397627.31
    97  = 4145162.075
    98  = 12996545.34
    99  = 464580.9394
    100 = 21032.94718
    101 = 25425029.74
    103 = 5769697.265
    104 = 5656888.156
    111 = 236325.294
    112 = 861150.6358
}

foreach ($row in $rowUpdates.Keys) {
    $newValue = $rowUpdates[$row]
    $rangeAddress = "J" + $row + ":AS" + $row
    $ws.Range($rangeAddress).Value = $newValue
}
